# Switching to Summer time
# Shift all timestamps in column A by +3 days (rows 2-97), and update
# the production values in column B for rows 2-56 to the new dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp (column A, rows 2 through 97) forward by 3 days.
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $cell.Value2 + 3
}

# New production values (MW) for column B, rows 2 through 56.
$newValues = @(
    1491, 1546, 1625, 1744, 1868, 1975, 2052, 2027, 2099, 2153,
    2213, 2250, 2350, 2359, 2371, 2387, 2446, 2453, 2439, 2467,
    2481, 2486, 2474, 2452, 2430, 2405, 2453, 2466, 2417, 2389,
    2363, 2300, 2299, 2319, 2324, 2342, 2342, 2339, 2350, 2334,
    2320, 2327, 2312, 2327, 2385, 2398, 2425, 2417, 2447, 2427,
    2421, 2411, 2422, 2432, 2465
)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 2).Value = $val
    $row++
}
